# Split the three run-on "Programa" / "Bibliografia" paragraphs into
# multiple numbered sentences separated by manual line breaks (<w:br/>),
# matching the target OOXML diff.

$d = $word.ActiveDocument

# --- Portuguese "Programa" paragraph -----------------------------------
$old1 = "1. Modelos PERT/COM2. Programação Linear Inteira; 2.1. Modelamento de problemas de PLI. 2.2 Algoritmo de ramificação e avaliação progressiva (branchand-bound).3. Programação Dinâmica3. Métodos Heurísticos; 3.1. Algoritmos Genéticos; 3.2 Recozimento Simulado; 3.3 Aplicação em problemas de otimização.4. Modelos e Técnicas de Previsão"
$new1 = "1. Modelos PERT/COM^l2. Programação Linear Inteira; 2.1. Modelamento de problemas de PLI. 2.2 Algoritmo de ramificação e avaliação progressiva (branchand-bound).^l3. Programação Dinâmica^l3. Métodos Heurísticos; 3.1. Algoritmos Genéticos; 3.2 Recozimento Simulado; 3.3 Aplicação em problemas de otimização.^l4. Modelos e Técnicas de Previsão"

$r1 = $d.Content
$r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- English "Programa" paragraph (italic) ------------------------------
$old2 = "1. PERT / COM models2. Whole Linear Programming; 2.1. Modeling of PLI problems. 2.2 Branch-bound algorithm.3. Dynamic Programming3. Heuristic methods; 3.1. Genetic Algorithms; 3.2 Simulated annealing; 3.3 Application in optimization problems.4. Forecasting Models and Techniques"
$new2 = "1. PERT / COM models^l2. Whole Linear Programming; 2.1. Modeling of PLI problems. 2.2 Branch-bound algorithm.^l3. Dynamic Programming^l3. Heuristic methods; 3.1. Genetic Algorithms; 3.2 Simulated annealing; 3.3 Application in optimization problems.^l4. Forecasting Models and Techniques"

$r2 = $d.Content
$r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Bibliografia paragraph ----------------------------------------------
$old3 = "1. HILLIER, F.S., LIEBERMAN, G.J., “Introdução à Pesquisa Operacional”, 8ªed., Editora McGraw-Hill, 2006.2. LACHTERMACHER, G., “Pesquisa Operacional na Tomada de Decisão (modelagem em Excel)”, 4ª ed., Editora Campus, 2009.3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., “An Introduction to Management Science” 9ª ed., South-Western College Publishing, 2000.4. PIZZOLATO, N. D. e GANDOLPHO, A. A. “Técnicas de Otimização”, LTC Editora, 2009.5. TAHA, H. A ., “Pesquisa Operacional”, 8ª ed., Pearson/Prentice Hall, 2008."
$new3 = "1. HILLIER, F.S., LIEBERMAN, G.J., “Introdução à Pesquisa Operacional”, 8ªed., Editora McGraw-Hill, 2006.^l2. LACHTERMACHER, G., “Pesquisa Operacional na Tomada de Decisão (modelagem em Excel)”, 4ª ed., Editora Campus, 2009.^l3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., “An Introduction to Management Science” 9ª ed., South-Western College Publishing, 2000.^l4. PIZZOLATO, N. D. e GANDOLPHO, A. A. “Técnicas de Otimização”, LTC Editora, 2009.^l5. TAHA, H. A ., “Pesquisa Operacional”, 8ª ed., Pearson/Prentice Hall, 2008."

$r3 = $d.Content
$r3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

Write-Host "Done"
